$wb = $excel.ActiveWorkbook

# Sheet "Chan1" - update rows 2 and 3 with new computed values
$ws1 = $wb.Worksheets.Item("Chan1")
$ws1.Range("A2").Value = 994
$ws1.Range("B2").Value = 331
$ws1.Range("C2").Value = 55
$ws1.Range("D2").Value = 608
$ws1.Range("E2").Value = 116
$ws1.Range("F2").Value = 24
$ws1.Range("G2").Value = 48
$ws1.Range("H2").Value = 44

$ws1.Range("A3").Value = 1319
$ws1.Range("B3").Value = 364
$ws1.Range("C3").Value = 345
$ws1.Range("D3").Value = 610
$ws1.Range("E3").Value = 251
$ws1.Range("F3").Value = 45
$ws1.Range("G3").Value = 139
$ws1.Range("H3").Value = 67

# Sheet "Chan2" - update rows 2 and 3 with new computed values
$ws2 = $wb.Worksheets.Item("Chan2")
$ws2.Range("A2").Value = 1267
$ws2.Range("B2").Value = 322
$ws2.Range("C2").Value = 234
$ws2.Range("D2").Value = 711
$ws2.Range("E2").Value = 140
$ws2.Range("F2").Value = 12
$ws2.Range("G2").Value = 88
$ws2.Range("H2").Value = 40

$ws2.Range("A3").Value = 1655
$ws2.Range("B3").Value = 342
$ws2.Range("C3").Value = 75
$ws2.Range("D3").Value = 1238
$ws2.Range("E3").Value = 180
$ws2.Range("F3").Value = 35
$ws2.Range("G3").Value = 56
$ws2.Range("H3").Value = 89
